$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.963.83'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '1.635.08'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("E4").Value = '  +0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '212.05'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.523'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  -0.01%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '23.43'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.07%  '
$ws.Range("E10").Value = '  -0.21%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0882'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").Value = '1.866.15'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = '1.632.45'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("E14").Value = '  -0.55%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.564'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.97%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '65.81'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '27.961.80'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '231.66'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.32%  '
$scratch = $ws.Range("Z1")
$scratch.Formula = '="0.0"&UNICHAR(8323)&"0725"'
$txt = $scratch.Text
$ws.Range("D19").Value = $txt
$scratch.Clear()
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("E21").Value = '  -0.03%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '10.42'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -6.59%  '
$ws.Range("E23").Value = '  -0.68%  '
$ws.Range("E24").Value = '  -0.29%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '155.05'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("E26").Value = '  +0.11%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '15.65'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.65%  '
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -0.51%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.40'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.65%  '
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("D34").Value = '1.406.03'
$ws.Range("E34").Value = '  -1.14%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.57'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.05'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +13.20%  '
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("E38").Value = '  +1.78%  '
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("E40").Value = '  -3.03%  '
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E46").Value = '  -0.26%  '
$ws.Range("D47").Value = '1.775.95'
$ws.Range("E47").Value = '  -0.45%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '88.25'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.77%  '
$scratch = $ws.Range("Z1")
$scratch.Formula = '="0.0"&UNICHAR(8326)&"0105"'
$txt = $scratch.Text
$ws.Range("D49").Value = $txt
$scratch.Clear()
$ws.Range("E49").Value = '  -0.89%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0999'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("E51").Value = '  -0.34%  '
